$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.965.32"
$ws.Range("E2").Value = "  +2.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.653.22"
$ws.Range("E3").Value = "  +2.76%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.89"
$ws.Range("E5").Value = "  +1.44%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +2.24%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +2.61%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.69%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.19"
$ws.Range("E10").Value = "  +4.92%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  +3.71%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.885.75"
$ws.Range("E12").Value = "  +2.80%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.653.28"
$ws.Range("E13").Value = "  +2.86%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +2.21%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.23"
$ws.Range("E16").Value = "  +2.89%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.963.21"
$ws.Range("E17").Value = "  +2.14%  "

# Row 18 - BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.48"
$ws.Range("E18").Value = "  +1.95%  "

# Row 19 and Row 20 - ShibaInu and Chainlink swap places
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.76"
$ws.Range("E19").Value = "  +2.05%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0733"
$ws.Range("E20").Value = "  +1.35%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.07%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +3.45%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +3.65%  "

# Row 24 - Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +1.17%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.92"
$ws.Range("E25").Value = "  -1.14%  "

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.12"
$ws.Range("E26").Value = "  +2.14%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.88%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +2.71%  "

# Row 29 - BinanceUSD
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.00%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.46%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.97%  "

# Row 32 - Maker
$ws.Range("D32").Value = "1.551.12"
$ws.Range("E32").Value = "  +3.89%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +2.68%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +9.87%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.13%  "

# Row 37 - ImmutableX
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.583"
$ws.Range("E37").Value = "  +3.37%  "

# Row 38 - ARBITRUM
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.901"
$ws.Range("E38").Value = "  +9.59%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +2.63%  "

# Row 40 - FraxShare
$ws.Range("E40").Value = "  +3.26%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.04%  "

# Row 42 - Aave
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.41"
$ws.Range("E42").Value = "  +9.09%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +2.67%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.794.35"
$ws.Range("E44").Value = "  +2.76%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.774"
$ws.Range("E45").Value = "  +1.59%  "

# Row 46 - WEMIXToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.936"
$ws.Range("E46").Value = "  +0.79%  "

# Row 47 - Quant
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.81"
$ws.Range("E47").Value = "  +0.42%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  +2.49%  "

# Row 49 - Algorand
$ws.Range("E49").Value = "  +3.03%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +0.87%  "

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  +2.19%  "
